# "fixed dates" - correct the two meeting dates in the timeline.
$d = $word.ActiveDocument

# First meeting: "November 7th" -> "November 17"
$d.Content.Find.Execute("First meeting (November 7th)", $true, $false, $false, $false, $false, $true, 1, $false, "First meeting (November 17)", 2)

# Third meeting: "December 2" -> "December 6"
$d.Content.Find.Execute("Third Meeting (December 2):", $true, $false, $false, $false, $false, $true, 1, $false, "Third Meeting (December 6):", 2)
